# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# completed handback: the Overview status moves from "Ready for handoff"
# to "Handed back: in sync with en-US", and both locale sheets (zh-cn,
# de-de) get their "Latest Target File" / "Latest Handback File" columns
# populated; de-de additionally gets a fresh "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$targetMdName = "d1d7db9c-142b-4746-a084-33bc9f3a0d22.md"
$targetMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8f4abdf5aaa564cc4a4838d879590fa450e095e/e2e/d1d7db9c-142b-4746-a084-33bc9f3a0d22.md"

# ---------------------------------------------------------------------
# Overview sheet: status text for both locale columns (E2, F2)
# ---------------------------------------------------------------------
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# widen the status columns to fit the longer text
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet: Latest Target File (I2) + Latest Handback File (J2)
# ---------------------------------------------------------------------
$zhcn.Range("I2").Value = $targetMdName
$zhcn.Range("I2").Style = "HyperLink"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $targetMdUrl, "", "", $targetMdName) | Out-Null

$zhcn.Range("J2").Value = "d1d7db9c-142b-4746-a084-33bc9f3a0d22.366bb1c57808d300d691e46ee4af2550ab0c8818.zh-cn.xlf"

# Latest Handback DateTime (K2) - zh-cn handback recorded
$zhcn.Range("K2").Value = "2016-08-25 15:02:50"

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet: Latest Target File (I2) + Latest Handback File (J2)
# ---------------------------------------------------------------------
$dede.Range("I2").Value = $targetMdName
$dede.Range("I2").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("I2"), $targetMdUrl, "", "", $targetMdName) | Out-Null

$dede.Range("J2").Value = "d1d7db9c-142b-4746-a084-33bc9f3a0d22.366bb1c57808d300d691e46ee4af2550ab0c8818.de-de.xlf"

# Latest Handback DateTime (K2) - de-de handback recorded (a few seconds after zh-cn)
$dede.Range("K2").Value = "2016-08-25 15:02:57"

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15
